$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 updates (target cluster ECs) ----
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,7).Value = 4.674406333333334
$ws.Cells.Item(2,8).Value = 14.023219
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,13).Value = 0.158112
$ws.Cells.Item(2,14).Value = 0.474336
$ws.Cells.Item(2,15).Value = 0.01629742442243621
$ws.Cells.Item(2,16).Value = 0.01629742442243621
$ws.Cells.Item(2,17).Value = 0.739079734176
$ws.Cells.Item(2,18).Value = 6.651717607584001
$ws.Cells.Item(2,19).Value = 0.01629742442243621
$ws.Cells.Item(2,20).Value = 0.01629742442243621

# ---- Row 3 updates (target cluster FAPs) ----
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,7).Value = 4.674406333333334
$ws.Cells.Item(3,8).Value = 14.023219
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,13).Value = 5.072121333333333
$ws.Cells.Item(3,14).Value = 15.216364
$ws.Cells.Item(3,15).Value = 0.5228098695318912
$ws.Cells.Item(3,16).Value = 0.5228098695318912
$ws.Cells.Item(3,17).Value = 23.70915608396844
$ws.Cells.Item(3,18).Value = 213.382404755716
$ws.Cells.Item(3,19).Value = 0.5228098695318912
$ws.Cells.Item(3,20).Value = 0.5228098695318912

# ---- Row 4: new row, target cluster M1 (new string) ----
$ws.Cells.Item(4,1).Value = "sCs"
$ws.Cells.Item(4,2).Value = "Bdnf"
$ws.Cells.Item(4,3).Value = "Ntrk2"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 4.674406333333334
$ws.Cells.Item(4,8).Value = 14.023219
$ws.Cells.Item(4,9).Value = 1
$ws.Cells.Item(4,10).Value = 1
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.010494
$ws.Cells.Item(4,14).Value = 0.031482
$ws.Cells.Item(4,15).Value = 0.001081671042609325
$ws.Cells.Item(4,16).Value = 0.001081671042609325
$ws.Cells.Item(4,17).Value = 0.049053220062
$ws.Cells.Item(4,18).Value = 0.441478980558
$ws.Cells.Item(4,19).Value = 0.001081671042609325
$ws.Cells.Item(4,20).Value = 0.001081671042609325

# ---- Row 5: new row, moved/updated old "sCs" target cluster row ----
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Bdnf"
$ws.Cells.Item(5,3).Value = "Ntrk2"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 4.674406333333334
$ws.Cells.Item(5,8).Value = 14.023219
$ws.Cells.Item(5,9).Value = 1
$ws.Cells.Item(5,10).Value = 1
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 4.460928333333333
$ws.Cells.Item(5,14).Value = 13.382785
$ws.Cells.Item(5,15).Value = 0.4598110350030633
$ws.Cells.Item(5,16).Value = 0.4598110350030633
$ws.Cells.Item(5,17).Value = 20.85219165387945
$ws.Cells.Item(5,18).Value = 187.669724884915
$ws.Cells.Item(5,19).Value = 0.4598110350030633
$ws.Cells.Item(5,20).Value = 0.4598110350030633
